# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-28 11:12:53
#
# This updates the "Recorded By" values (column G) on the "Session Analysis
# Results" sheet. A handful of rows had their recorded-by name lists
# reordered (the "System" entry moved from the front of the comma-separated
# list to the end) when the report was re-synced from the source system.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows where "System, backup@backdoor.com, system" -> "system, System, backup@backdoor.com"
$rowsBackdoor = @(2, 28, 54)

# Rows where "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
$rowsDnasr = @(3, 6, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 26, 29, 32, 36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 52, 55, 58, 62, 63, 64, 65, 66, 67, 69, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86, 90, 92, 93, 94, 96, 99, 101, 109, 110, 111, 112, 116, 118, 119, 120, 122, 125, 127, 135, 136, 137, 138, 142, 144, 145, 146, 148, 151, 153)

# Rows where "System, admin@admin.com" -> "admin@admin.com, System"
$rowsAdmin = @(7, 33, 59)

$oldBackdoor = "System, backup@backdoor.com, system"
$newBackdoor = "system, System, backup@backdoor.com"

$oldDnasr = "System, dnasr281@gmail.com"
$newDnasr = "dnasr281@gmail.com, System"

$oldAdmin = "System, admin@admin.com"
$newAdmin = "admin@admin.com, System"

foreach ($r in $rowsBackdoor) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldBackdoor) {
        $cell.Value = $newBackdoor
    } else {
        Write-Host "Row $r unexpected value: $($cell.Value2)"
    }
}

foreach ($r in $rowsDnasr) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldDnasr) {
        $cell.Value = $newDnasr
    } else {
        Write-Host "Row $r unexpected value: $($cell.Value2)"
    }
}

foreach ($r in $rowsAdmin) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldAdmin) {
        $cell.Value = $newAdmin
    } else {
        Write-Host "Row $r unexpected value: $($cell.Value2)"
    }
}

Write-Host "Done."
